# Applies the "Organização dos arquivos para rodar novo fluxo 5 nós" edit:
#  - BESS!B2: bus-node id 4 -> shared string "650_sec" (style cleared)
#  - Load!B2:B41: numeric bus-node ids (e.g. 650001) -> "<3-digit-prefix>_sec" strings
#  - Load!K2:K41: terminal codes (e.g. 234) -> truncated by one digit (234 -> 23)
#  - Public_Ilumination!B2:B14: numeric bus-node ids (e.g. 650) -> "<id>_sec" strings
#  - Selection / active-sheet bookkeeping matching the new focus (Load sheet active)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# BESS sheet: B2 goes from a plain number (4, bold style) to the shared
# string "650_sec". Clear the existing formatting first so the cell reverts
# to the default style, matching the target (style attribute dropped).
# ---------------------------------------------------------------------------
$wsBESS = $wb.Worksheets.Item("BESS")
$wsBESS.Range("B2").Clear()
$wsBESS.Range("B2").Value = "650_sec"

# ---------------------------------------------------------------------------
# Load sheet: 40 data rows. Column B holds bus-node ids like 650001, 632002,
# etc. -- the new value is the 3-digit substation prefix plus "_sec". Column
# K holds 3-digit terminal codes like 234 -- the new value drops the last
# digit (234 -> 23, 124 -> 12, 134 -> 13).
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Load")
for ($r = 2; $r -le 41; $r++) {
    $busCell = $wsLoad.Cells.Item($r, 2)
    $busId = [int]$busCell.Value()
    $prefix = [math]::Floor($busId / 1000)
    $busCell.Value = "$($prefix)_sec"

    $termCell = $wsLoad.Cells.Item($r, 11)
    $termVal = [int]$termCell.Value()
    $termCell.Value = [math]::Floor($termVal / 10)
}

# ---------------------------------------------------------------------------
# Public_Ilumination sheet: 13 data rows. Column B already holds the bare
# 3-digit bus-node id (650, 632, ...) -- append "_sec".
# ---------------------------------------------------------------------------
$wsPI = $wb.Worksheets.Item("Public_Ilumination")
for ($r = 2; $r -le 14; $r++) {
    $busCell = $wsPI.Cells.Item($r, 2)
    $busId = [int]$busCell.Value()
    $busCell.Value = "$($busId)_sec"
}

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the workbook focus moves from
# Generator to Load, and each sheet keeps a fresh selection anchor.
# ---------------------------------------------------------------------------
$wsBESS.Activate()
$wsBESS.Range("C6").Select() | Out-Null

$wsGen = $wb.Worksheets.Item("Generator")
$wsGen.Activate()
$wsGen.Range("G20").Select() | Out-Null

$wsPI.Activate()
$wsPI.Range("C24").Select() | Out-Null

$wsLoad.Activate()
$wsLoad.Range("O14").Select() | Out-Null
